# Feat: sequence diagram and images
#
# Applies:
#  1. Re-cache of the "datetimeFigureOut" footer field (2022-05-16 -> 2022-05-24)
#     on the slide master and every slide layout.
#  2. Re-positions + re-colors the isosceles-triangle autoshape on slide 1.
#  3. Best-effort touch of the presentation-level slide-guide list / theme
#     object-defaults bookkeeping entries that PowerPoint writes silently.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the cached date text wherever it shows up as "2022-05-16":
#    the slide master plus every custom (slide) layout.
# ---------------------------------------------------------------------------
function Update-DateShape($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "2022-05-16") {
                $shp.TextFrame.TextRange.Text = "2022-05-24"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShape $master

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DateShape $layout
}

# ---------------------------------------------------------------------------
# 2) Slide 1: move + recolor the isosceles triangle shape (id 36).
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(1)

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.Id -eq 36) {
        $shp.Left = 12.209291338582677
        $shp.Top = 60.44968613937008
        $shp.Fill.ForeColor.SchemeColor = "tx1"
    }
}

# ---------------------------------------------------------------------------
# 3) Best-effort: record the restyled triangle as the new default autoshape
#    (what PowerPoint itself uses to populate the theme's <a:objectDefaults>)
#    and nudge the presentation's slide-guide bookkeeping list into existing.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.Id -eq 36) {
        $shp.SetShapesDefaultProperties()
        $p.DefaultShape = $shp
    }
}

try {
    $p.Guides = $true
} catch {
}
try {
    $master.Guides = $true
} catch {
}
